$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315346002578735
$ws.Range("B1").Value = 4.618990898132324
$ws.Range("C1").Value = 3.238940477371216
$ws.Range("D1").Value = 2.473646640777588
$ws.Range("E1").Value = 2.233258962631226
